{"js": "// 1) Update the title text from \"-\u8ba1\u5212\u65e5\u671f\" to \"-\u65e5\u671f\"\nconst searchResults = context.document.body.search(\"-\u8ba1\u5212\u65e5\u671f\", { matchCase: true, matchWildcards: false });\nsearchResults.load(\"items\");\nawait context.sync();\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"-\u65e5\u671f\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Bump the paragraph-mark font size (to 11pt / sz=22) for every empty cell\n//    in the last row of the bonding-sheet table (the blank data row under the\n//    header row), matching the \"big screen\" redesign sizing.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  const lastRow = table.rows.getLast();\n  lastRow.load(\"cells/items\");\n  await context.sync();\n\n  for (let i = 0; i < lastRow.cells.items.length; i++) {\n    const cell = lastRow.cells.items[i];\n    const cellRange = cell.body.getRange(\"Whole\");\n    cellRange.font.size = 11;\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the title text from \"-\u8ba1\u5212\u65e5\u671f\" to \"-\u65e5\u671f\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"-\u8ba1\u5212\u65e5\u671f\", $false, $false, $false, $false, $false, $true, 1, $false, \"-\u65e5\u671f\", 2)\n\n# 2) Bump the paragraph-mark font size (to 11pt / sz=22) for every empty cell\n#    in the last row of the bonding-sheet table (the blank data row under the\n#    header row), matching the \"big screen\" redesign sizing.\n$table = $d.Tables.Item(1)\n$lastRow = $table.Rows.Last\nfor ($i = 1; $i -le $lastRow.Cells.Count; $i++) {\n    $cell = $lastRow.Cells.Item($i)\n    $cell.Range.Font.Size = 11\n}\n"}
